# Add new "Player Info" worksheet as the first sheet in the workbook.
$wb = $excel.ActiveWorkbook

$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

# Data row
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4324"
$playerInfo.Range("B2").Value = "Taijul Islam"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Rename the MATCH_CARD_LINK column to MATCH_CODE and replace the full
# scorecard URL with just the numeric match code, on both the
# "ODI Batting" and "ODI Bowling" sheets. Look the sheets up by name
# since indices shift once the new sheet has been inserted.
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

$matchCodes = @("3708","3784","3936","3938","4357","4358","4416","4418","4420","4616","4627","4628","4711","4713","4717")

# ODI Batting: MATCH_CARD_LINK is column D (4)
$odiBatting.Cells.Item(1, 4).Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $odiBatting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

# ODI Bowling: MATCH_CARD_LINK is column B (2)
$odiBowling.Cells.Item(1, 2).Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $odiBowling.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

Write-Host "Edit complete"
